# Continue setting up the Config workbook: add a new Constants entry for
# fetching the Finance PSTST credentials asset, and register the four new
# Orchestrator assets (environment, login URL, email-to field, automation
# environment) used by the "Update Physical Asset History" process.

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets = $wb.Worksheets.Item("Assets")

# --- Assets sheet: the new rows all share the same Orchestrator folder, so
# fill that column first -----------------------------------------------------
$wsAssets.Range("C2:C5").Value = "UITS"

# --- Constants sheet: new row 18 -------------------------------------------------
$wsConstants.Range("A18").Value = "FinancePhysicalAssetHistory_Credentials"
$wsConstants.Range("B18").Value = "FinancePhysicalAssetHistory_Credentials"
$wsConstants.Range("C18").Value = "the name of the asset that holds the credentials for the Finance PSTST environment"

# --- Assets sheet: new rows 2-5 ---------------------------------------------------
$wsAssets.Range("A2").Value = "FinancePhysicalAssetHistory_Environment"
$wsAssets.Range("B2").Value = "FinancePhysicalAssetHistory_Environment"
$wsAssets.Range("D2").Value = "The Finance PSTST Environment"

$wsAssets.Range("A3").Value = "FinancePhysicalAssetHistory_LoginURL"
$wsAssets.Range("B3").Value = "FinancePhysicalAssetHistory_LoginURL"
$wsAssets.Range("D3").Value = "The URL to log into Finance PSTST"

$wsAssets.Range("A4").Value = "FinancePhysicalAssetHistory_EmailToField"
$wsAssets.Range("B4").Value = "FinancePhysicalAssetHistory_EmailToField"
$wsAssets.Range("D4").Value = "A list of semicolon-separated emails to send the start & end emails to"

$wsAssets.Range("A5").Value = "FinancePhysicalAssetHistory_AutomationEnvironment"
$wsAssets.Range("B5").Value = "FinancePhysicalAssetHistory_AutomationEnvironment"
$wsAssets.Range("D5").Value = "The environment of the automation, either UAT or Production"

# --- Leave the UI state reflecting where the author's cursor ended up ------------
$wsConstants.Range("A19").Select() | Out-Null
$wsAssets.Activate() | Out-Null
$wsAssets.Range("D5").Select() | Out-Null
